# Update "想去人数" (interest count) figures in column F across the four
# worksheets, reflecting the refreshed scrape output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 27
$ws.Range("F6").Value  = 1055
$ws.Range("F9").Value  = 522
$ws.Range("F11").Value = 414
$ws.Range("F12").Value = 147
$ws.Range("F13").Value = 1256
$ws.Range("F14").Value = 1162
$ws.Range("F15").Value = 1339
$ws.Range("F16").Value = 256
$ws.Range("F17").Value = 1495
$ws.Range("F20").Value = 283
$ws.Range("F23").Value = 997
$ws.Range("F25").Value = 773
$ws.Range("F27").Value = 913
$ws.Range("F28").Value = 172842
$ws.Range("F29").Value = 890
$ws.Range("F32").Value = 867
$ws.Range("F34").Value = 1498
$ws.Range("F35").Value = 73
$ws.Range("F36").Value = 757
$ws.Range("F38").Value = 745

# Sheet: 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 102
$ws.Range("F6").Value  = 132
$ws.Range("F11").Value = 1306
$ws.Range("F13").Value = 2398
$ws.Range("F17").Value = 158
$ws.Range("F18").Value = 25
$ws.Range("F19").Value = 54
$ws.Range("F22").Value = 407
$ws.Range("F25").Value = 246
$ws.Range("F29").Value = 217
$ws.Range("F31").Value = 29
$ws.Range("F37").Value = 26

# Sheet: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value  = 2636
$ws.Range("F7").Value  = 4379
$ws.Range("F8").Value  = 100
$ws.Range("F10").Value = 468
$ws.Range("F11").Value = 496
$ws.Range("F12").Value = 358
$ws.Range("F13").Value = 407
$ws.Range("F14").Value = 135

# Sheet: 全部类型 (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value  = 4379
$ws.Range("F6").Value  = 100
$ws.Range("F7").Value  = 468
$ws.Range("F8").Value  = 496
$ws.Range("F9").Value  = 27
$ws.Range("F10").Value = 407
$ws.Range("F11").Value = 135
$ws.Range("F12").Value = 132
$ws.Range("F14").Value = 1055
$ws.Range("F17").Value = 1306
$ws.Range("F18").Value = 522
$ws.Range("F20").Value = 414
$ws.Range("F21").Value = 147
$ws.Range("F22").Value = 2398
$ws.Range("F24").Value = 1162
$ws.Range("F25").Value = 1339
$ws.Range("F27").Value = 158
$ws.Range("F28").Value = 54
$ws.Range("F29").Value = 1495
$ws.Range("F31").Value = 283
$ws.Range("F32").Value = 407
$ws.Range("F33").Value = 997
$ws.Range("F34").Value = 773
$ws.Range("F36").Value = 913
$ws.Range("F37").Value = 246
$ws.Range("F38").Value = 890
$ws.Range("F40").Value = 867
$ws.Range("F43").Value = 1498
$ws.Range("F44").Value = 73
$ws.Range("F47").Value = 757
$ws.Range("F48").Value = 745
